$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 14:04"

# Row 5
$ws.Cells.Item(5, 2).Value = 260117
$ws.Cells.Item(5, 3).Value = 3262
$ws.Cells.Item(5, 4).Value = 168408
$ws.Cells.Item(5, 5).Value = 65410
$ws.Cells.Item(5, 7).Value = 229
$ws.Cells.Item(5, 8).Value = 26299

# Row 11
$ws.Cells.Item(11, 2).Value = 136519
$ws.Cells.Item(11, 3).Value = 826
$ws.Cells.Item(11, 5).Value = 71904
$ws.Cells.Item(11, 7).Value = 77
$ws.Cells.Item(11, 8).Value = 9265

# Row 22
$ws.Cells.Item(22, 2).Value = 30207
$ws.Cells.Item(22, 3).Value = 81
$ws.Cells.Item(22, 5).Value = 2497

# Row 32
$ws.Cells.Item(32, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(32, 2).Value = 16793
$ws.Cells.Item(32, 3).Value = 553
$ws.Cells.Item(32, 4).Value = 3837
$ws.Cells.Item(32, 5).Value = 12782
$ws.Cells.Item(32, 6).Value = 1
$ws.Cells.Item(32, 7).Value = 9
$ws.Cells.Item(32, 8).Value = 174

# Row 33
$ws.Cells.Item(33, 1).Value = "Israel"
$ws.Cells.Item(33, 2).Value = 16409
$ws.Cells.Item(33, 3).Value = 28
$ws.Cells.Item(33, 4).Value = 11007
$ws.Cells.Item(33, 5).Value = 5157
$ws.Cells.Item(33, 6).Value = 77
$ws.Cells.Item(33, 7).Value = 5
$ws.Cells.Item(33, 8).Value = 245

# Row 84
$ws.Cells.Item(84, 2).Value = 1586
$ws.Cells.Item(84, 3).Value = 14
$ws.Cells.Item(84, 4).Value = 1099
$ws.Cells.Item(84, 5).Value = 397
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 90

# Row 86
$ws.Cells.Item(86, 2).Value = 1551
$ws.Cells.Item(86, 3).Value = 59
$ws.Cells.Item(86, 4).Value = 611
$ws.Cells.Item(86, 5).Value = 927

# Row 120
$ws.Cells.Item(120, 1).Value = "Malta"
$ws.Cells.Item(120, 2).Value = 489
$ws.Cells.Item(120, 3).Value = 3
$ws.Cells.Item(120, 4).Value = 419
$ws.Cells.Item(120, 5).Value = 65
$ws.Cells.Item(120, 8).Value = 5

# Row 121
$ws.Cells.Item(121, 1).Value = "Jamaica"
$ws.Cells.Item(121, 2).Value = 488
$ws.Cells.Item(121, 3).Value = 10
$ws.Cells.Item(121, 4).Value = 58
$ws.Cells.Item(121, 5).Value = 421
$ws.Cells.Item(121, 8).Value = 9

# Row 140
$ws.Cells.Item(140, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(140, 2).Value = 208
$ws.Cells.Item(140, 3).Value = 21
$ws.Cells.Item(140, 4).Value = 4
$ws.Cells.Item(140, 5).Value = 199
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 1
$ws.Cells.Item(140, 8).Value = 5

# Row 141
$ws.Cells.Item(141, 1).Value = "Etiopia"
$ws.Cells.Item(141, 2).Value = 194
$ws.Cells.Item(141, 3).Value = 3
$ws.Cells.Item(141, 4).Value = 95
$ws.Cells.Item(141, 5).Value = 95
$ws.Cells.Item(141, 8).Value = 4

# Row 142
$ws.Cells.Item(142, 1).Value = "Madagascar"
$ws.Cells.Item(142, 2).Value = 193
$ws.Cells.Item(142, 4).Value = 101
$ws.Cells.Item(142, 5).Value = 92
$ws.Cells.Item(142, 6).Value = 1
$ws.Cells.Item(142, 8).Value = 0

# Row 143
$ws.Cells.Item(143, 1).Value = "Liberia"
$ws.Cells.Item(143, 2).Value = 189
$ws.Cells.Item(143, 4).Value = 79
$ws.Cells.Item(143, 5).Value = 90
$ws.Cells.Item(143, 8).Value = 20

# Row 158
$ws.Cells.Item(158, 2).Value = 102
$ws.Cells.Item(158, 3).Value = 1
$ws.Cells.Item(158, 4).Value = 30
$ws.Cells.Item(158, 5).Value = 72

# Row 205
$ws.Cells.Item(205, 1).Value = "Seychelles"
$ws.Cells.Item(205, 4).Value = 8
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 8).Value = 0

# Row 206
$ws.Cells.Item(206, 1).Value = "Montserrat"
$ws.Cells.Item(206, 4).Value = 7
$ws.Cells.Item(206, 6).Value = 1
$ws.Cells.Item(206, 8).Value = 1
